$wb = $excel.ActiveWorkbook

# "Ready for handoff" -> "Handback transform failed" everywhere it appears
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667
$wsZhCn.Range("P3").Value = "Handback file name: vgzotgm3.ib0 is different with handoff file name: af038034-0e16-43af-a25c-7229d97a965c.f3e00e0aab1a0aa1d1fd9e73e18e7671fd0fa6e5.zh-cn."

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
$wsDeDe.Range("P3").Value = "Handback file name: vgzotgm3.ib0 is different with handoff file name: af038034-0e16-43af-a25c-7229d97a965c.f3e00e0aab1a0aa1d1fd9e73e18e7671fd0fa6e5.de-de."
